$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.748.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.20%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.458.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.77%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.77%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.458.82'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.75%  '

# Row 8
$ws.Range("E8").Value = '  -0.17%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.66%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.122'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.76%  '

# Row 12
$ws.Range("E12").Value = '  -7.88%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.040.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.82%  '

# Row 14
$ws.Range("E14").Value = '  -10.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -10.02%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.687.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.32%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.422.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.57%  '

# Row 18
$ws.Range("E18").Value = '  -2.33%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -10.59%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.56%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.73%  '

# Row 23
$ws.Range("E23").Value = '  -10.34%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.94%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.603.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.59%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000107'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.86%  '

# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.68%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.57%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.462.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.60%  '

# Row 33
$ws.Range("E33").Value = '  -0.02%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.147'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.06%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '173.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.03%  '

# Row 37
$ws.Range("E37").Value = '  -13.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.95'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.51%  '

# Row 40
$ws.Range("E40").Value = '  -12.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0779'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.64%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.55%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '

# Row 45
$ws.Range("E45").Value = '  -14.23%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.93%  '

# Row 48
$ws.Range("E48").Value = '  -2.31%  '

# Row 49
$ws.Range("E49").Value = '  -7.61%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.220.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.96%  '
